# "Generate Report for handback"
#
# The localization-status workbook lists 3 source files (ffff594e..., ffffffd8...,
# 6228f1a5...) across 3 sheets (Overview, zh-cn, de-de). This edit reflects that the
# 6228f1a5 file has now been "Handed back: in sync with en-US" (it used to be only
# "Ready for handoff", with no handback info yet). As a result all rows are re-sorted
# alphabetically by source file name on every sheet, and the 6228f1a5 row gets its
# full handback data (Latest Target File / Latest Handback File / Latest Handback
# DateTime) filled in, matching the shape of the other, already-handed-back rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

$ov.Cells.Item(2, 2).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3, 2).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(4, 2).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/ffff594e8794-256f-480c-acc1-8f630d18b6d6.md", "", "", "ffff594e8794-256f-480c-acc1-8f630d18b6d6.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/ffffffd8ab7988-9f04-4e6b-a8d7-0abd9cf45f44.md", "", "", "ffffffd8ab7988-9f04-4e6b-a8d7-0abd9cf45f44.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File | Latest Handoff
# Datetime | Latest Target File | Latest Handback File | Latest Handback DateTime |
# Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

# Row 2: 6228f1a5 (now handed back)
$zh.Cells.Item(2, 2).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(2, 4).Value = "2016-01-25 14:06:50"
$zh.Cells.Item(2, 7).Value = "2016-01-25 14:07:29"
$zh.Cells.Item(2, 8).Value = "Include"

# Row 3: ffff594e8794 (same content it had, just a different row now)
$zh.Cells.Item(3, 2).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(3, 4).Value = "2016-01-25 14:04:38"
$zh.Cells.Item(3, 7).Value = "2016-01-25 14:05:24"
$zh.Cells.Item(3, 8).Value = "Include"

# Row 4: ffffffd8ab79
$zh.Cells.Item(4, 2).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(4, 4).Value = "2016-01-25 14:04:38"
$zh.Cells.Item(4, 7).Value = "2016-01-25 14:05:24"
$zh.Cells.Item(4, 8).Value = "Include"

# Row 2 hyperlinks/text
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/911f765a8d9ab825dd6583396eb18839f9d990b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.e42cba792d8a499b05907ba00b3273ff373b857b.zh-cn.xlf", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.e42cba792d8a499b05907ba00b3273ff373b857b.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/911f765a8d9ab825dd6583396eb18839f9d990b3/e2e/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/911f765a8d9ab825dd6583396eb18839f9d990b3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.e42cba792d8a499b05907ba00b3273ff373b857b.zh-cn.xlf", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.e42cba792d8a499b05907ba00b3273ff373b857b.zh-cn.xlf") | Out-Null

# Row 3 hyperlinks/text
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/ffff594e8794-256f-480c-acc1-8f630d18b6d6.md", "", "", "ffff594e8794-256f-480c-acc1-8f630d18b6d6.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd31d5f684cb1211ade3a31a86bab3224db93055/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/957e3db8384645b59405d7802299ae166227f3d0/e2e/19610356-b2fd-4781-a25d-de0dee7cda64.md", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fdf4772c664414c402c1a5d7435cb1fe9ccb96c6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf") | Out-Null

# Row 4 hyperlinks/text
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/ffffffd8ab7988-9f04-4e6b-a8d7-0abd9cf45f44.md", "", "", "ffffffd8ab7988-9f04-4e6b-a8d7-0abd9cf45f44.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd31d5f684cb1211ade3a31a86bab3224db93055/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/957e3db8384645b59405d7802299ae166227f3d0/e2e/19610356-b2fd-4781-a25d-de0dee7cda64.md", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fdf4772c664414c402c1a5d7435cb1fe9ccb96c6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

# Row 2: 6228f1a5 (now handed back)
$de.Cells.Item(2, 2).Value = "Handed back: in sync with en-US"
$de.Cells.Item(2, 4).Value = "2016-01-25 14:07:06"
$de.Cells.Item(2, 7).Value = "2016-01-25 14:07:47"
$de.Cells.Item(2, 8).Value = "Include"

# Row 3: ffff594e8794
$de.Cells.Item(3, 2).Value = "Handed back: in sync with en-US"
$de.Cells.Item(3, 4).Value = "2016-01-25 14:04:47"
$de.Cells.Item(3, 7).Value = "2016-01-25 14:05:41"
$de.Cells.Item(3, 8).Value = "Include"

# Row 4: ffffffd8ab79
$de.Cells.Item(4, 2).Value = "Handed back: in sync with en-US"
$de.Cells.Item(4, 4).Value = "2016-01-25 14:04:47"
$de.Cells.Item(4, 7).Value = "2016-01-25 14:05:41"
$de.Cells.Item(4, 8).Value = "Include"

# Row 2 hyperlinks/text
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/30f24c61afa1ef26054a05a83c57c28737e4d4c7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.e42cba792d8a499b05907ba00b3273ff373b857b.de-de.xlf", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.e42cba792d8a499b05907ba00b3273ff373b857b.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/30f24c61afa1ef26054a05a83c57c28737e4d4c7/e2e/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/30f24c61afa1ef26054a05a83c57c28737e4d4c7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.e42cba792d8a499b05907ba00b3273ff373b857b.de-de.xlf", "", "", "6228f1a5-a9ba-4cd6-bd3c-9aa77302112b.e42cba792d8a499b05907ba00b3273ff373b857b.de-de.xlf") | Out-Null

# Row 3 hyperlinks/text
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/ffff594e8794-256f-480c-acc1-8f630d18b6d6.md", "", "", "ffff594e8794-256f-480c-acc1-8f630d18b6d6.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a31e88b108b13a934237e8a8d9af7d34759e0e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8ad7951d9f7f9b6aaad0024b10aa301fa216e81f/e2e/19610356-b2fd-4781-a25d-de0dee7cda64.md", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d7e9110608fc42292e405f83625e12c19b10d801/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf") | Out-Null

# Row 4 hyperlinks/text
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/e2e/ffffffd8ab7988-9f04-4e6b-a8d7-0abd9cf45f44.md", "", "", "ffffffd8ab7988-9f04-4e6b-a8d7-0abd9cf45f44.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a31e88b108b13a934237e8a8d9af7d34759e0e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8ad7951d9f7f9b6aaad0024b10aa301fa216e81f/e2e/19610356-b2fd-4781-a25d-de0dee7cda64.md", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d7e9110608fc42292e405f83625e12c19b10d801/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf", "", "", "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ec433aa6bbda04e21d299ae2947f7bcee4d40887/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Report generated for handback."
